# Admin adding special dates, updated logic on displaying meals and menu
#
# The weekly menu sheet ("List2") is refreshed:
#   - the "From"/"To" dates move to the following week
#   - the food (Jidlo) descriptions are replaced with the new Czech menu
#     item names ("... / English name" placeholders for translation)
#   - the soup (Polivka) descriptions are replaced with the new Czech
#     soup names, now listed in the same top-to-bottom order as the days
#   - row heights for the wrapped-text rows are nudged to their new
#     auto-fit sizes
#   - the last-used selection is moved to K4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Special dates for this week's menu (column A = Od/From, B = Do/To) ---
$ws.Range("A2").Value = 45908
$ws.Range("B2").Value = 45912

# --- Food names (column C) ---
$ws.Range("C2").Value = "Něměcký Řízek s bramborovou kaší / English name"
$ws.Range("C3").Value = "Český řízek s bramborovou kaší / English name"
$ws.Range("C4").Value = "Maďarský guláš / English name"
$ws.Range("C5").Value = "Pečený candát na másle / English name"
$ws.Range("C6").Value = "Srbská pljeskavica / English name"

# --- Soup names (column F), now aligned with the day order in rows 2..6 ---
$ws.Range("F2").Value = "Německá polévka"
$ws.Range("F3").Value = "Česká polévka"
$ws.Range("F4").Value = "Maďarská polévka"
$ws.Range("F5").Value = "Švédská polévka"
$ws.Range("F6").Value = "Srbská polévka"

# --- Row heights re-settle to their new auto-fit values after the text edit ---
$ws.Rows.Item(2).RowHeight = 94.2
$ws.Rows.Item(5).RowHeight = 55.95

# --- Leave the selection where the admin last clicked ---
[void]$ws.Range("K4").Select()
